$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '302.61'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-5.16%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '35.11'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-2.73%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.044'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-1.71%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07963'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-3.11%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.938'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-9.94%'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '7.744'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-3.28%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.954'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '5.48%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9224'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-0.35%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1221'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '21.56%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1837'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-3.02%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.09397'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '2.20%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03527'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-1.94%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09853'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-0.67%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001405'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-1.98%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005844'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '3.21%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.497'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '0.96%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.049'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-2.13%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3447'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2.14%'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-3.01%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.044'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-0.37%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04494'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-2.07%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001215'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-2.47%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004850'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '2.36%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001249'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-3.90%'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-6.81%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01914'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '-4.77%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04747'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-4.83%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007525'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-2.87%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.009555'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '22.39%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1326'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-5.31%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002108'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '0.62%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.01113'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-7.11%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006287'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-2.61%'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000749'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.06%'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-31.34%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002098'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.06%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0001999'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.06%'
